$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Reorder work-experience sections: move the "Analytics Supervisor -
#    GSD&M" block (heading + 4 body paragraphs) so it comes right after the
#    "Data Products Manager - Helm/Murmuration" block instead of before it.
# ---------------------------------------------------------------------------

$gsdmStart = $null
$gsdmEndIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Analytics Supervisor - GSD*") {
        $gsdmStart = $i
        break
    }
}

$gsdmEndIndex = $gsdmStart + 4   # heading + 4 bullet/body paragraphs = 5 total

$startPar = $d.Paragraphs($gsdmStart)
$endPar = $d.Paragraphs($gsdmEndIndex)
$moveRange = $d.Range($startPar.Range.Start, $endPar.Range.End)
$moveRange.Cut() | Out-Null

# Find the new location: the paragraph ending in "... by 57%" (last bullet of
# the "Data Products Manager" section), which - after the cut above - sits
# right before what is now the "Senior Analyst - Myers Research" heading.
$dpmEndIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Senior Analyst - Myers Research*") {
        $dpmEndIndex = $i - 1
        break
    }
}

$target = $d.Paragraphs($dpmEndIndex)
$insertPoint = $d.Range($target.Range.End, $target.Range.End)
$insertPoint.Paste() | Out-Null

# The paste loses the "Heading 3" style on the first pasted paragraph
# (a quirk of this COM-interop's Paste implementation) - restore it.
$pastedHeading = $d.Paragraphs($dpmEndIndex + 1)
$pastedHeading.Style = "Heading 3"

# ---------------------------------------------------------------------------
# 2) Work-experience bullet: neutralize the language, with "50M" becoming its
#    own bold, colored run (matching the styling of the other stat call-outs
#    in this bullet, e.g. "23%"/"64%"). Do this FIRST, while "50M" is still
#    unique in the document, so the subsequent formatting Find() unambiguously
#    lands on this occurrence rather than one introduced by the plain-text
#    replacements below.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial machine learning",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial machine learning",
    2) | Out-Null

$statRange = $d.Content
$statRange.Find.Execute("50M", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$statRange.Font.Bold = 1
$statRange.Font.Color = 5258796   # 0x2C3E50 stored BGR, matches the other stat runs

# ---------------------------------------------------------------------------
# 3) Professional summary: neutralize the demographic-coding-error language
#    (plain text, no special run formatting).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Key-projects "Impact" line for the Geospatial Demographic Classification
#    System project.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2) | Out-Null
